# Auto-generated Excel COM-interop script applying the Sagittarius_Profits data refresh
# to sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW (per scheduled-runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 792.5
$ws.Range("J17").Value = 817.26086
$ws.Range("L17").Value = 2451.78258
$ws.Range("N17").Value = -2787.78258
$ws.Range("H51").Value = 3740
$ws.Range("I51").Value = 3740
$ws.Range("K51").Value = 3740
$ws.Range("M51").Value = -3256
$ws.Range("H76").Value = 2950
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 2900
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 2900
$ws.Range("M76").Value = -2685
$ws.Range("N76").Value = -3530
$ws.Range("H79").Value = 2950
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 2900
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 2900
$ws.Range("M79").Value = -1908
$ws.Range("N79").Value = -5084
$ws.Range("H127").Value = 2011.5454
$ws.Range("I127").Value = 2190.5715
$ws.Range("J127").Value = 1698.25
$ws.Range("K127").Value = 6571.7145
$ws.Range("L127").Value = 5094.75
$ws.Range("M127").Value = -1611.7145
$ws.Range("N127").Value = -15014.75
$ws.Range("H137").Value = 1543.2222
$ws.Range("I137").Value = 1559.6
$ws.Range("K137").Value = 4678.799999999999
$ws.Range("M137").Value = -2128.799999999999
$ws.Range("H138").Value = 5081.775
$ws.Range("J138").Value = 5633.3335
$ws.Range("L138").Value = 16900.0005
$ws.Range("N138").Value = -27180.0005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5137253.5
$ws.Range("I32").Value = 5004200.5
$ws.Range("K32").Value = 5004200.5
$ws.Range("M32").Value = -5003913.5
$ws.Range("H110").Value = 4500
$ws.Range("I110").Value = 4500
$ws.Range("K110").Value = 4500
$ws.Range("M110").Value = -2455
$ws.Range("H111").Value = 27663.334
$ws.Range("J111").Value = 27663.334
$ws.Range("L111").Value = 27663.334
$ws.Range("N111").Value = -35843.334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2300.2727
$ws.Range("I20").Value = 1651.6666
$ws.Range("K20").Value = 1651.6666
$ws.Range("M20").Value = -1404.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 11572.5
$ws.Range("I88").Value = 12000
$ws.Range("J88").Value = 11511.429
$ws.Range("K88").Value = 12000
$ws.Range("L88").Value = 11511.429
$ws.Range("M88").Value = -11594
$ws.Range("N88").Value = -12323.429
$ws.Range("H91").Value = 11572.5
$ws.Range("I91").Value = 12000
$ws.Range("J91").Value = 11511.429
$ws.Range("K91").Value = 12000
$ws.Range("L91").Value = 11511.429
$ws.Range("M91").Value = -10596
$ws.Range("N91").Value = -14319.429
$ws.Range("H107").Value = 330.5
$ws.Range("I107").Value = 295.75
$ws.Range("K107").Value = 295.75
$ws.Range("M107").Value = 1624.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.90476
$ws.Range("I2").Value = 20.333334
$ws.Range("J2").Value = 168.66667
$ws.Range("K2").Value = 122.000004
$ws.Range("L2").Value = 1012.00002
$ws.Range("M2").Value = -9.000004000000004
$ws.Range("N2").Value = -1238.00002
$ws.Range("H4").Value = 9166771
$ws.Range("I4").Value = 12222306
$ws.Range("K4").Value = 36666918
$ws.Range("M4").Value = -36666806
$ws.Range("H7").Value = 805.625
$ws.Range("I7").Value = 250
$ws.Range("J7").Value = 1731.6666
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 5194.9998
$ws.Range("M7").Value = -638
$ws.Range("N7").Value = -5418.9998
$ws.Range("H12").Value = 152.9
$ws.Range("I12").Value = 217.2
$ws.Range("J12").Value = 88.59999999999999
$ws.Range("K12").Value = 651.5999999999999
$ws.Range("L12").Value = 265.8
$ws.Range("M12").Value = -478.5999999999999
$ws.Range("N12").Value = -611.8
$ws.Range("H33").Value = 1032.3334
$ws.Range("J33").Value = 999.5
$ws.Range("L33").Value = 5997
$ws.Range("N33").Value = -6563
$ws.Range("H109").Value = 1141
$ws.Range("I109").Value = 201.25
$ws.Range("J109").Value = 4900
$ws.Range("K109").Value = 603.75
$ws.Range("L109").Value = 14700
$ws.Range("M109").Value = 436.25
$ws.Range("N109").Value = -16780
$ws.Range("H113").Value = 1126.3125
$ws.Range("J113").Value = 1178.3
$ws.Range("L113").Value = 3534.9
$ws.Range("N113").Value = -7874.9
$ws.Range("H115").Value = 3959.8
$ws.Range("I115").Value = 3899
$ws.Range("J115").Value = 3975
$ws.Range("K115").Value = 11697
$ws.Range("L115").Value = 11925
$ws.Range("M115").Value = -10522
$ws.Range("N115").Value = -14275
$ws.Range("H140").Value = 10937.667
$ws.Range("I140").Value = 938.75
$ws.Range("J140").Value = 15937.125
$ws.Range("K140").Value = 2816.25
$ws.Range("L140").Value = 47811.375
$ws.Range("M140").Value = 2363.75
$ws.Range("N140").Value = -58171.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 134167.33
$ws.Range("I7").Value = 199999.5
$ws.Range("J7").Value = 2503
$ws.Range("K7").Value = 199999.5
$ws.Range("L7").Value = 2503
$ws.Range("M7").Value = -199887.5
$ws.Range("N7").Value = -2727
$ws.Range("H8").Value = 134167.33
$ws.Range("I8").Value = 199999.5
$ws.Range("J8").Value = 2503
$ws.Range("K8").Value = 199999.5
$ws.Range("L8").Value = 2503
$ws.Range("M8").Value = -199860.5
$ws.Range("N8").Value = -2781
$ws.Range("H9").Value = 935
$ws.Range("I9").Value = 98.5
$ws.Range("J9").Value = 2608
$ws.Range("K9").Value = 98.5
$ws.Range("L9").Value = 2608
$ws.Range("M9").Value = 71.5
$ws.Range("N9").Value = -2948
$ws.Range("H11").Value = 6667501.5
$ws.Range("I11").Value = 10000000
$ws.Range("J11").Value = 2504
$ws.Range("K11").Value = 10000000
$ws.Range("L11").Value = 2504
$ws.Range("M11").Value = -9999861
$ws.Range("N11").Value = -2782
$ws.Range("H12").Value = 2504
$ws.Range("J12").Value = 2504
$ws.Range("L12").Value = 2504
$ws.Range("N12").Value = -2784
$ws.Range("H14").Value = 500666.66
$ws.Range("I14").Value = 500666.66
$ws.Range("K14").Value = 500666.66
$ws.Range("M14").Value = -500498.66
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H19").Value = 885
$ws.Range("I19").Value = 541.5833
$ws.Range("J19").Value = 5006
$ws.Range("K19").Value = 541.5833
$ws.Range("L19").Value = 5006
$ws.Range("M19").Value = -253.5833
$ws.Range("N19").Value = -5582
$ws.Range("H52").Value = 20000
$ws.Range("I52").Value = 20000
$ws.Range("K52").Value = 20000
$ws.Range("M52").Value = -19741
$ws.Range("H63").Value = 38556.5
$ws.Range("J63").Value = 38556.5
$ws.Range("L63").Value = 38556.5
$ws.Range("N63").Value = -39928.5
$ws.Range("H66").Value = 38556.5
$ws.Range("J66").Value = 38556.5
$ws.Range("L66").Value = 115669.5
$ws.Range("N66").Value = -122533.5
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
$ws.Range("H132").Value = 19989
$ws.Range("J132").Value = 19988
$ws.Range("L132").Value = 59964
$ws.Range("N132").Value = -65024
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 325.66666
$ws.Range("I16").Value = 325.66666
$ws.Range("K16").Value = 325.66666
$ws.Range("M16").Value = -155.66666
$ws.Range("H46").Value = 2203.6365
$ws.Range("J46").Value = 2259.25
$ws.Range("L46").Value = 2259.25
$ws.Range("N46").Value = -2635.25
$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15450
$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16560
$ws.Range("H122").Value = 6904.3335
$ws.Range("I122").Value = 5181.727
$ws.Range("J122").Value = 8799.200000000001
$ws.Range("K122").Value = 15545.181
$ws.Range("L122").Value = 26397.6
$ws.Range("M122").Value = -13095.181
$ws.Range("N122").Value = -31297.6
